# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with newer case counts.
#
# A handful of countries received updated totals (column B, "Casos totales")
# that push them past one or two neighbours in the table, which is sorted
# descending by column B. We update those countries' stats in place (by
# name, wherever they currently live) and then re-sort the whole data
# block so the ranking reflects the new numbers - exactly like a refresh
# of the live feed that backs this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A4:H219")
$countryCol = $ws.Range("A4:A219")

function Set-CountryStats {
    param(
        [string]$Country,
        [double]$CasosTotales,
        [double]$NuevosCasos,
        [double]$CasosActivos,
        [double]$Recuperados,
        [double]$CasosCriticos,
        [double]$MuertesHoy,
        [double]$Muertes
    )

    $cell = $countryCol.Find($Country)
    $row = $cell.Row

    $ws.Cells.Item($row, 2).Value = $CasosTotales
    $ws.Cells.Item($row, 3).Value = $NuevosCasos
    $ws.Cells.Item($row, 4).Value = $CasosActivos
    $ws.Cells.Item($row, 5).Value = $Recuperados
    $ws.Cells.Item($row, 6).Value = $CasosCriticos
    $ws.Cells.Item($row, 7).Value = $MuertesHoy
    $ws.Cells.Item($row, 8).Value = $Muertes
}

Set-CountryStats "Estados Unidos"        6007864  7499  3315042 2508975 0 194 183847
Set-CountryStats "Alemania"               239608   608   214233   16019 0   4   9356
Set-CountryStats "Republica Dominicana"    92964   407    64347   26987 0  17   1630
Set-CountryStats "Portugal"                56673   399    41357   13507 0   2   1809
Set-CountryStats "Azerbaiyan"              35844   137    33364    1956 0   2    524
Set-CountryStats "Moldavia"                35546   564    24156   10413 0  10    977
Set-CountryStats "Noruega"                 10524    20     9348     912 0   0    264
Set-CountryStats "Albania"                  9083   156     4791    4026 0   3    266
Set-CountryStats "Tayikistan"               8449    36     7246    1135 0   1     68
Set-CountryStats "Namibia"                  6712   281     2776    3876 0   1     60
Set-CountryStats "Jordania"                 1801    45     1364     422 0   0     15
Set-CountryStats "Birmania"                  602    22      345     251 0   0      6

# Re-sort the table by "Casos totales" (column B) descending, same as the
# live sheet does on every refresh, so the newly-updated countries land in
# their correct rank (this reorders rows without touching any other data).
$sortKey = $ws.Range("B4")
$dataRange.Sort($sortKey, 2, $null, $null, 1, $null, 1, 1)

# Update the "last refreshed" timestamp shown at the top of the sheet.
$ws.Range("A1").Value = "Datos actualizados a 27 de Agosto de 2020 a las 17:11"
